$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 5.041666666666667
$ws.Range("C2").Value = 6

$ws.Range("B3").Value = 5.041666666666667
$ws.Range("C3").Value = 6

$ws.Range("B4").Value = 3.875
$ws.Range("C4").Value = 4

$ws.Range("B5").Value = 8.208333333333334
$ws.Range("C5").Value = 9

$ws.Range("B6").Value = 18.625
$ws.Range("C6").Value = 19

$ws.Range("B7").Value = 3.458333333333333
$ws.Range("C7").Value = 4
